$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Jogos"

# Keep the pivot cache's worksheet source pointing at the renamed sheet.
$ws2 = $wb.Worksheets.Item(2)
$pt = $ws2.PivotTables(1)
try {
    $pt.SourceData = "Jogos!A1:R1048576"
} catch {
}

$ws.Range("E8").Select()
